$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7213929295539856
$ws.Range("B1").Value = 1.37474524974823
$ws.Range("C1").Value = 4.31864595413208
$ws.Range("D1").Value = 1.921665191650391
$ws.Range("E1").Value = 0.9638192653656006
